# add code for cham cong function
# Update last_edited_time stamps (D3, D4, D5, D7, D13) and the
# "Chi tieu" / "Luy ke" figures for Thang 7 (row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2024-07-06T13:10:00.000Z"

$ws.Range("D3").Value = $newTimestamp
$ws.Range("D4").Value = $newTimestamp
$ws.Range("D5").Value = $newTimestamp
$ws.Range("D7").Value = $newTimestamp
$ws.Range("D13").Value = $newTimestamp

$ws.Range("W5").Value = 9335000
$ws.Range("AA5").Value = -4035000
